# Update Facebook and Twitter social media analytics data (rows 4-10)
# Columns C:O = msg_count_twitter group, P:AB = msg_count_twitter_engage group,
# AC:AO = msg_count_facebook group (see row 1 merged headers / row 2 stat labels).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 794859
$ws.Range("D4").Value = 1290.4
$ws.Range("E4").Value = 2353.9
$ws.Range("H4").Value = 510
$ws.Range("I4").Value = 1756.2
$ws.Range("J4").Value = 23112
$ws.Range("L4").Value = 1967.5
$ws.Range("M4").Value = 404
$ws.Range("N4").Value = 65.59999999999999
$ws.Range("O4").Value = 0.4
$ws.Range("P4").Value = 289499
$ws.Range("Q4").Value = 470
$ws.Range("R4").Value = 1595.4
$ws.Range("U4").Value = 63
$ws.Range("V4").Value = 348.2
$ws.Range("W4").Value = 20340
$ws.Range("Y4").Value = 704.4
$ws.Range("Z4").Value = 411
$ws.Range("AA4").Value = 66.7
$ws.Range("AB4").Value = 0.4
$ws.Range("AC4").Value = 399102
$ws.Range("AD4").Value = 647.9
$ws.Range("AE4").Value = 820.8
$ws.Range("AH4").Value = 389
$ws.Range("AI4").Value = 1074.5
$ws.Range("AJ4").Value = 8295
$ws.Range("AL4").Value = 919.6
$ws.Range("AM4").Value = 434
$ws.Range("AN4").Value = 70.5
$ws.Range("AO4").Value = 0.8

# Row 5
$ws.Range("C5").Value = 888413
$ws.Range("D5").Value = 920.6
$ws.Range("E5").Value = 2920
$ws.Range("H5").Value = 257
$ws.Range("I5").Value = 967
$ws.Range("J5").Value = 55081
$ws.Range("L5").Value = 1322
$ws.Range("M5").Value = 672
$ws.Range("N5").Value = 69.59999999999999
$ws.Range("P5").Value = 270025
$ws.Range("Q5").Value = 279.8
$ws.Range("R5").Value = 895.7
$ws.Range("U5").Value = 37
$ws.Range("V5").Value = 209
$ws.Range("W5").Value = 10903
$ws.Range("Y5").Value = 395.9
$ws.Range("Z5").Value = 682
$ws.Range("AA5").Value = 70.7
$ws.Range("AC5").Value = 433055
$ws.Range("AD5").Value = 448.8
$ws.Range("AE5").Value = 599.7
$ws.Range("AH5").Value = 254
$ws.Range("AI5").Value = 652
$ws.Range("AJ5").Value = 5325
$ws.Range("AL5").Value = 623.1
$ws.Range("AM5").Value = 695
$ws.Range("AN5").Value = 72
$ws.Range("AO5").Value = 0.9

# Row 6
$ws.Range("C6").Value = 165551
$ws.Range("D6").Value = 2365
$ws.Range("E6").Value = 2812.9
$ws.Range("G6").Value = 326.8
$ws.Range("H6").Value = 1632
$ws.Range("I6").Value = 3697.8
$ws.Range("J6").Value = 14651
$ws.Range("L6").Value = 2956.3
$ws.Range("M6").Value = 56
$ws.Range("N6").Value = 80
$ws.Range("P6").Value = 99578
$ws.Range("Q6").Value = 1422.5
$ws.Range("R6").Value = 6031.9
$ws.Range("T6").Value = 56.2
$ws.Range("U6").Value = 241.5
$ws.Range("V6").Value = 659.8
$ws.Range("W6").Value = 48717
$ws.Range("Y6").Value = 1778.2
$ws.Range("Z6").Value = 56
$ws.Range("AA6").Value = 80
$ws.Range("AC6").Value = 45694
$ws.Range("AD6").Value = 652.8
$ws.Range("AE6").Value = 661
$ws.Range("AG6").Value = 0.5
$ws.Range("AH6").Value = 436.5
$ws.Range("AI6").Value = 1141.5
$ws.Range("AJ6").Value = 2204
$ws.Range("AL6").Value = 878.7
$ws.Range("AM6").Value = 52
$ws.Range("AN6").Value = 74.3
$ws.Range("AO6").Value = 1.1

# Row 7
$ws.Range("C7").Value = 196201
$ws.Range("D7").Value = 1096.1
$ws.Range("E7").Value = 2726.4
$ws.Range("H7").Value = 31
$ws.Range("I7").Value = 1012
$ws.Range("J7").Value = 24791
$ws.Range("L7").Value = 2065.3
$ws.Range("O7").Value = -0.5
$ws.Range("P7").Value = 77644
$ws.Range("Q7").Value = 433.8
$ws.Range("R7").Value = 1448.7
$ws.Range("U7").Value = 5
$ws.Range("V7").Value = 136.5
$ws.Range("W7").Value = 10080
$ws.Range("Y7").Value = 817.3
$ws.Range("AC7").Value = 95551
$ws.Range("AD7").Value = 533.8
$ws.Range("AE7").Value = 876.5
$ws.Range("AH7").Value = 25
$ws.Range("AI7").Value = 787
$ws.Range("AJ7").Value = 6193
$ws.Range("AL7").Value = 1038.6
$ws.Range("AM7").Value = 92
$ws.Range("AN7").Value = 51.4
$ws.Range("AO7").Value = -0.6

# Row 8
$ws.Range("C8").Value = 729459
$ws.Range("D8").Value = 509.4
$ws.Range("E8").Value = 1796.5
$ws.Range("I8").Value = 427
$ws.Range("J8").Value = 52263
$ws.Range("L8").Value = 1010.3
$ws.Range("M8").Value = 722
$ws.Range("N8").Value = 50.4
$ws.Range("P8").Value = 309148
$ws.Range("Q8").Value = 215.9
$ws.Range("R8").Value = 1286.3
$ws.Range("V8").Value = 72
$ws.Range("W8").Value = 24618
$ws.Range("Y8").Value = 415.5
$ws.Range("Z8").Value = 744
$ws.Range("AA8").Value = 52
$ws.Range("AC8").Value = 469182
$ws.Range("AD8").Value = 327.6
$ws.Range("AE8").Value = 550.4
$ws.Range("AH8").Value = 48.5
$ws.Range("AI8").Value = 468.8
$ws.Range("AJ8").Value = 5226
$ws.Range("AL8").Value = 563.9
$ws.Range("AM8").Value = 832
$ws.Range("AN8").Value = 58.1
$ws.Range("AO8").Value = -0.1

# Row 9
$ws.Range("C9").Value = 22909
$ws.Range("D9").Value = 1041.3
$ws.Range("E9").Value = 1313.5
$ws.Range("H9").Value = 257.5
$ws.Range("I9").Value = 2184.2
$ws.Range("J9").Value = 4285
$ws.Range("L9").Value = 1636.4
$ws.Range("P9").Value = 4431
$ws.Range("Q9").Value = 201.4
$ws.Range("R9").Value = 352.3
$ws.Range("U9").Value = 21
$ws.Range("V9").Value = 233
$ws.Range("W9").Value = 1437
$ws.Range("Y9").Value = 316.5
$ws.Range("AC9").Value = 6924
$ws.Range("AD9").Value = 314.7
$ws.Range("AE9").Value = 535.8
$ws.Range("AI9").Value = 455.8
$ws.Range("AJ9").Value = 1845
$ws.Range("AL9").Value = 865.5
$ws.Range("AM9").Value = 8
$ws.Range("AN9").Value = 36.4
$ws.Range("AO9").Value = -1.7

# Row 10
$ws.Range("C10").Value = 25556
$ws.Range("D10").Value = 464.7
$ws.Range("E10").Value = 1305.2
$ws.Range("I10").Value = 130
$ws.Range("J10").Value = 6144
$ws.Range("L10").Value = 1161.6
$ws.Range("P10").Value = 3929
$ws.Range("Q10").Value = 71.40000000000001
$ws.Range("R10").Value = 198.9
$ws.Range("V10").Value = 8.5
$ws.Range("W10").Value = 1010
$ws.Range("Y10").Value = 187.1
$ws.Range("Z10").Value = 21
$ws.Range("AA10").Value = 38.2
$ws.Range("AB10").Value = -0.5
$ws.Range("AC10").Value = 15178
$ws.Range("AD10").Value = 276
$ws.Range("AE10").Value = 583.9
$ws.Range("AH10").Value = 6
$ws.Range("AI10").Value = 240
$ws.Range("AJ10").Value = 2668
$ws.Range("AL10").Value = 523.4
$ws.Range("AM10").Value = 29
$ws.Range("AN10").Value = 52.7
$ws.Range("AO10").Value = -0.5
